$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 310.66666
$ws.Range("I2").Value = 234.5
$ws.Range("J2").Value = 371.6
$ws.Range("K2").Value = 234.5
$ws.Range("L2").Value = 371.6
$ws.Range("M2").Value = -121.5
$ws.Range("N2").Value = -597.6
$ws.Range("H43").Value = 1550.3636
$ws.Range("I43").Value = 950
$ws.Range("J43").Value = 1610.4
$ws.Range("K43").Value = 950
$ws.Range("L43").Value = 1610.4
$ws.Range("M43").Value = -881
$ws.Range("N43").Value = -1748.4
$ws.Range("H100").Value = 2824.6
$ws.Range("I100").Value = 2763.7856
$ws.Range("J100").Value = 2966.5
$ws.Range("K100").Value = 2763.7856
$ws.Range("L100").Value = 2966.5
$ws.Range("M100").Value = -2222.7856
$ws.Range("N100").Value = -4048.5
$ws.Range("H112").Value = 4334.125
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 4500.2075
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 13500.6225
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -15716.6225
$ws.Range("H129").Value = 1198.5106
$ws.Range("I129").Value = 546.4
$ws.Range("J129").Value = 1374.7567
$ws.Range("K129").Value = 1639.2
$ws.Range("L129").Value = 4124.2701
$ws.Range("M129").Value = 3360.8
$ws.Range("N129").Value = -14124.2701
$ws.Range("H138").Value = 1963.24
$ws.Range("I138").Value = 599.93335
$ws.Range("J138").Value = 2203.8235
$ws.Range("K138").Value = 1799.80005
$ws.Range("L138").Value = 6611.470499999999
$ws.Range("M138").Value = 3340.19995
$ws.Range("N138").Value = -16891.4705

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 389343.2
$ws.Range("I32").Value = 422246.94
$ws.Range("J32").Value = 27402
$ws.Range("K32").Value = 422246.94
$ws.Range("L32").Value = 27402
$ws.Range("M32").Value = -421959.94
$ws.Range("N32").Value = -27976
$ws.Range("H45").Value = 3443.8572
$ws.Range("I45").Value = 2544.8572
$ws.Range("J45").Value = 4342.857
$ws.Range("K45").Value = 2544.8572
$ws.Range("L45").Value = 4342.857
$ws.Range("M45").Value = -2167.8572
$ws.Range("N45").Value = -5096.857
$ws.Range("H64").Value = 38000
$ws.Range("J64").Value = 38000
$ws.Range("L64").Value = 38000
$ws.Range("N64").Value = -38496
$ws.Range("H67").Value = 38000
$ws.Range("J67").Value = 38000
$ws.Range("L67").Value = 38000
$ws.Range("N67").Value = -39716
$ws.Range("H122").Value = 127262.5
$ws.Range("I122").Value = 201580
$ws.Range("K122").Value = 604740
$ws.Range("M122").Value = -602290

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H94").Value = 835.73914
$ws.Range("I94").Value = 650.6
$ws.Range("K94").Value = 650.6
$ws.Range("M94").Value = -199.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 284.66666
$ws.Range("I2").Value = 284.66666
$ws.Range("K2").Value = 284.66666
$ws.Range("M2").Value = -171.66666
$ws.Range("H3").Value = 1304
$ws.Range("I3").Value = 1304
$ws.Range("K3").Value = 1304
$ws.Range("M3").Value = -1191
$ws.Range("H5").Value = 1761.6666
$ws.Range("I5").Value = 368
$ws.Range("J5").Value = 2757.1428
$ws.Range("K5").Value = 368
$ws.Range("L5").Value = 2757.1428
$ws.Range("M5").Value = -256
$ws.Range("N5").Value = -2981.1428
$ws.Range("H11").Value = 1970.3334
$ws.Range("I11").Value = 305
$ws.Range("J11").Value = 2803
$ws.Range("K11").Value = 305
$ws.Range("L11").Value = 2803
$ws.Range("M11").Value = -165
$ws.Range("N11").Value = -3083
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H31").Value = 5988.0815
$ws.Range("I31").Value = 1670.1072
$ws.Range("K31").Value = 1670.1072
$ws.Range("M31").Value = -1375.1072
$ws.Range("H34").Value = 5988.0815
$ws.Range("I34").Value = 1670.1072
$ws.Range("K34").Value = 1670.1072
$ws.Range("M34").Value = -1468.1072
$ws.Range("H58").Value = 1279.7587
$ws.Range("I58").Value = 940.25
$ws.Range("K58").Value = 940.25
$ws.Range("M58").Value = -737.25
$ws.Range("H136").Value = 1279.7587
$ws.Range("I136").Value = 940.25
$ws.Range("K136").Value = 2820.75
$ws.Range("M136").Value = -270.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 587.5
$ws.Range("I5").Value = 528.5714
$ws.Range("K5").Value = 1585.7142
$ws.Range("M5").Value = -1473.7142
$ws.Range("H23").Value = 58823660
$ws.Range("I23").Value = 111.166664
$ws.Range("J23").Value = 90909230
$ws.Range("K23").Value = 333.499992
$ws.Range("L23").Value = 272727690
$ws.Range("M23").Value = -98.49999200000002
$ws.Range("N23").Value = -272728160
$ws.Range("H68").Value = 868
$ws.Range("J68").Value = 1000
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4622
$ws.Range("H71").Value = 868
$ws.Range("J71").Value = 1000
$ws.Range("L71").Value = 9000
$ws.Range("N71").Value = -17112
$ws.Range("H92").Value = 866
$ws.Range("J92").Value = 866
$ws.Range("L92").Value = 2598
$ws.Range("N92").Value = -5094
$ws.Range("H122").Value = 7889.5
$ws.Range("I122").Value = 340
$ws.Range("K122").Value = 3060
$ws.Range("M122").Value = -610
$ws.Range("H135").Value = 587.5
$ws.Range("I135").Value = 528.5714
$ws.Range("K135").Value = 4757.1426
$ws.Range("M135").Value = -2222.1426
$ws.Range("H138").Value = 7693.4165
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 7693.4165
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 23080.2495
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -33360.24950000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3174.353
$ws.Range("I122").Value = 2398.2222
$ws.Range("J122").Value = 4047.5
$ws.Range("K122").Value = 7194.6666
$ws.Range("L122").Value = 12142.5
$ws.Range("M122").Value = -4744.6666
$ws.Range("N122").Value = -17042.5
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1375.75
$ws.Range("I7").Value = 1045.8182
$ws.Range("K7").Value = 1045.8182
$ws.Range("M7").Value = -933.8181999999999
$ws.Range("H122").Value = 3322.6
$ws.Range("I122").Value = 3200.2
$ws.Range("J122").Value = 3567.4
$ws.Range("K122").Value = 9600.599999999999
$ws.Range("L122").Value = 10702.2
$ws.Range("M122").Value = -7150.599999999999
$ws.Range("N122").Value = -15602.2
$ws.Range("H126").Value = 1375.75
$ws.Range("I126").Value = 1045.8182
$ws.Range("K126").Value = 3137.4546
$ws.Range("M126").Value = -667.4546
$ws.Range("H132").Value = 4059.318
$ws.Range("I132").Value = 3491.9167
$ws.Range("K132").Value = 10475.7501
$ws.Range("M132").Value = -7945.750100000001
$ws.Range("H136").Value = 6668745.5
$ws.Range("I136").Value = 1874.9231
$ws.Range("J136").Value = 13891189
$ws.Range("K136").Value = 5624.7693
$ws.Range("L136").Value = 41673567
$ws.Range("M136").Value = -3074.7693
$ws.Range("N136").Value = -41678667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 800.4
$ws.Range("I107").Value = 667.3333
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2001.9999
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -81.99990000000003
$ws.Range("N107").Value = -6840
$ws.Range("H122").Value = 2859.875
$ws.Range("I122").Value = 1977
$ws.Range("J122").Value = 3154.1667
$ws.Range("K122").Value = 5931
$ws.Range("L122").Value = 9462.500100000001
$ws.Range("M122").Value = -3481
$ws.Range("N122").Value = -14362.5001
$ws.Range("H123").Value = 22691.5
$ws.Range("J123").Value = 22691.5
$ws.Range("L123").Value = 22691.5
$ws.Range("N123").Value = -32491.5
$ws.Range("H136").Value = 4398.4116
$ws.Range("I136").Value = 4524.364
$ws.Range("K136").Value = 13573.092
$ws.Range("M136").Value = -11023.092
